# Update a set of imputed numeric values in the KNN result data sheet.
# (commit: "Update Name of Algo" -- underlying algorithm output values changed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.352
$ws.Range("B4").Value = 4.914
$ws.Range("A6").Value = -22.178
$ws.Range("A7").Value = -20.978
$ws.Range("D7").Value = -8.007
$ws.Range("D8").Value = -8.648
$ws.Range("B9").Value = 6.322000000000001
$ws.Range("D10").Value = -8.071000000000002
$ws.Range("B12").Value = 5.609
$ws.Range("D13").Value = -8.260999999999999
$ws.Range("A16").Value = -21.24
$ws.Range("D16").Value = -8.430000000000001
$ws.Range("B17").Value = 5.179
$ws.Range("B18").Value = 5.790000000000001
$ws.Range("B19").Value = 7.491
$ws.Range("A20").Value = -21.814
$ws.Range("B20").Value = 4.996
$ws.Range("B26").Value = 6.115
$ws.Range("A28").Value = -21.604
$ws.Range("A29").Value = -21.5
$ws.Range("D30").Value = -7.37
$ws.Range("B31").Value = 6.630999999999998
$ws.Range("A32").Value = -21.418
$ws.Range("B39").Value = 7.145
$ws.Range("A40").Value = -21.021
$ws.Range("B40").Value = 7.042
$ws.Range("D40").Value = -8.301
$ws.Range("B41").Value = 6.251
$ws.Range("B42").Value = 6.101
$ws.Range("B43").Value = 6.786
$ws.Range("D44").Value = -7.510999999999998
$ws.Range("A46").Value = -21.489
$ws.Range("B47").Value = 6.146000000000001
$ws.Range("B48").Value = 6.258999999999999
$ws.Range("A51").Value = -21.557
$ws.Range("A52").Value = -21.646
$ws.Range("A57").Value = -21.806
$ws.Range("A59").Value = -22.257
$ws.Range("A62").Value = -21.696
$ws.Range("B63").Value = 5.583
$ws.Range("B64").Value = 6.259
$ws.Range("A66").Value = -21.527
$ws.Range("A73").Value = -21.099
$ws.Range("A74").Value = -20.458
$ws.Range("B76").Value = 5.792
$ws.Range("B81").Value = 5.468999999999999
$ws.Range("B89").Value = 5.641999999999999
$ws.Range("D89").Value = -8.347
$ws.Range("D91").Value = -7.634
$ws.Range("A92").Value = -21.566
$ws.Range("B94").Value = 5.923
$ws.Range("A100").Value = -21.961
